$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value2 = 168.4
$ws.Range("H17").Value2 = 2214.3572
$ws.Range("J17").Value2 = 2307.6924
$ws.Range("L17").Value2 = 6923.0772
$ws.Range("N17").Value2 = -7259.0772
$ws.Range("H132").Value2 = 1673.5094
$ws.Range("I132").Value2 = 1427.2291
$ws.Range("K132").Value2 = 4281.6873
$ws.Range("M132").Value2 = -1751.6873
$ws.Range("H138").Value2 = 6192.875
$ws.Range("I138").Value2 = 3732.5789
$ws.Range("J138").Value2 = 6799.961
$ws.Range("K138").Value2 = 11197.7367
$ws.Range("L138").Value2 = 20399.883
$ws.Range("M138").Value2 = -6057.736699999999
$ws.Range("N138").Value2 = -30679.883
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 5591.989
$ws.Range("I32").Value2 = 4210.2705
$ws.Range("K32").Value2 = 4210.2705
$ws.Range("M32").Value2 = -3923.2705
$ws.Range("H45").Value2 = 2811.7144
$ws.Range("I45").Value2 = 2378
$ws.Range("K45").Value2 = 2378
$ws.Range("M45").Value2 = -2001
$ws.Range("H61").Value2 = 5157
$ws.Range("I61").Value2 = 4952.7
$ws.Range("K61").Value2 = 4952.7
$ws.Range("M61").Value2 = -4740.7
$ws.Range("H74").Value2 = 2003.2931
$ws.Range("I74").Value2 = 2048.84
$ws.Range("K74").Value2 = 2048.84
$ws.Range("M74").Value2 = -1174.84
$ws.Range("H77").Value2 = 2003.2931
$ws.Range("I77").Value2 = 2048.84
$ws.Range("K77").Value2 = 10244.2
$ws.Range("M77").Value2 = -5876.200000000001
$ws.Range("H119").Value2 = 149250
$ws.Range("J119").Value2 = 149250
$ws.Range("L119").Value2 = 149250
$ws.Range("N119").Value2 = -158926
$ws.Range("H122").Value2 = 4110.8213
$ws.Range("I122").Value2 = 1983.5834
$ws.Range("K122").Value2 = 5950.7502
$ws.Range("M122").Value2 = -3500.7502
$ws.Range("H132").Value2 = 3749.4204
$ws.Range("I132").Value2 = 3057.2856
$ws.Range("K132").Value2 = 9171.856800000001
$ws.Range("M132").Value2 = -6641.856800000001
$ws.Range("H136").Value2 = 5157
$ws.Range("I136").Value2 = 4952.7
$ws.Range("K136").Value2 = 14858.1
$ws.Range("M136").Value2 = -12308.1
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value2 = 1180.7142
$ws.Range("I11").Value2 = 88.666664
$ws.Range("J11").Value2 = 1999.75
$ws.Range("K11").Value2 = 88.666664
$ws.Range("L11").Value2 = 1999.75
$ws.Range("M11").Value2 = 51.333336
$ws.Range("N11").Value2 = -2279.75
$ws.Range("H99").Value2 = 1666.6666
$ws.Range("I99").Value2 = 1500
$ws.Range("J99").Value2 = 2000
$ws.Range("K99").Value2 = 1500
$ws.Range("L99").Value2 = 2000
$ws.Range("M99").Value2 = -2
$ws.Range("N99").Value2 = -4996
$ws.Range("H105").Value2 = 170070
$ws.Range("I105").Value2 = 253152.5
$ws.Range("J105").Value2 = 3905
$ws.Range("K105").Value2 = 253152.5
$ws.Range("L105").Value2 = 3905
$ws.Range("M105").Value2 = -251405.5
$ws.Range("N105").Value2 = -7399
$ws.Range("H134").Value2 = 19732.133
$ws.Range("I134").Value2 = 2981.5173
$ws.Range("K134").Value2 = 8944.5519
$ws.Range("M134").Value2 = -6409.5519
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value2 = 603
$ws.Range("I22").Value2 = 415.36365
$ws.Range("J22").Value2 = 947
$ws.Range("K22").Value2 = 415.36365
$ws.Range("L22").Value2 = 947
$ws.Range("M22").Value2 = -65.36365000000001
$ws.Range("N22").Value2 = -1647
$ws.Range("H31").Value2 = 59152.777
$ws.Range("I31").Value2 = 1887.5454
$ws.Range("J31").Value2 = 149141
$ws.Range("K31").Value2 = 1887.5454
$ws.Range("L31").Value2 = 149141
$ws.Range("M31").Value2 = -1592.5454
$ws.Range("N31").Value2 = -149731
$ws.Range("H34").Value2 = 59152.777
$ws.Range("I34").Value2 = 1887.5454
$ws.Range("J34").Value2 = 149141
$ws.Range("K34").Value2 = 1887.5454
$ws.Range("L34").Value2 = 149141
$ws.Range("M34").Value2 = -1685.5454
$ws.Range("N34").Value2 = -149545
$ws.Range("H132").Value2 = 4679.905
$ws.Range("I132").Value2 = 3127
$ws.Range("J132").Value2 = 7785.7144
$ws.Range("K132").Value2 = 9381
$ws.Range("L132").Value2 = 23357.1432
$ws.Range("M132").Value2 = -6851
$ws.Range("N132").Value2 = -28417.1432
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value2 = 3149.5
$ws.Range("I69").Value2 = 899.6667
$ws.Range("K69").Value2 = 2699.0001
$ws.Range("M69").Value2 = -1888.0001
$ws.Range("H70").Value2 = 168599.5
$ws.Range("I70").Value2 = 168599.5
$ws.Range("K70").Value2 = 505798.5
$ws.Range("M70").Value2 = -505483.5
$ws.Range("H72").Value2 = 3149.5
$ws.Range("I72").Value2 = 899.6667
$ws.Range("K72").Value2 = 8097.0003
$ws.Range("M72").Value2 = -4041.0003
$ws.Range("H73").Value2 = 168599.5
$ws.Range("I73").Value2 = 168599.5
$ws.Range("K73").Value2 = 505798.5
$ws.Range("M73").Value2 = -504706.5
$ws.Range("H74").Value2 = 0
$ws.Range("I74").Value2 = 0
$ws.Range("K74").Value2 = 0
$ws.Range("M74").ClearContents() | Out-Null
$ws.Range("H77").Value2 = 0
$ws.Range("I77").Value2 = 0
$ws.Range("K77").Value2 = 0
$ws.Range("M77").ClearContents() | Out-Null
$ws.Range("H87").Value2 = 33098.516
$ws.Range("I87").Value2 = 25812.5
$ws.Range("J87").Value2 = 34999.22
$ws.Range("K87").Value2 = 77437.5
$ws.Range("L87").Value2 = 104997.66
$ws.Range("M87").Value2 = -76189.5
$ws.Range("N87").Value2 = -107493.66
$ws.Range("H90").Value2 = 33098.516
$ws.Range("I90").Value2 = 25812.5
$ws.Range("J90").Value2 = 34999.22
$ws.Range("K90").Value2 = 232312.5
$ws.Range("L90").Value2 = 314992.98
$ws.Range("M90").Value2 = -226072.5
$ws.Range("N90").Value2 = -327472.98
$ws.Range("H112").Value2 = 252824.75
$ws.Range("I112").Value2 = 501249.5
$ws.Range("J112").Value2 = 4400
$ws.Range("K112").Value2 = 1503748.5
$ws.Range("L112").Value2 = 13200
$ws.Range("M112").Value2 = -1502640.5
$ws.Range("N112").Value2 = -15416
$ws.Range("H137").Value2 = 6117.7
$ws.Range("J137").Value2 = 7474.4165
$ws.Range("L137").Value2 = 22423.2495
$ws.Range("N137").Value2 = -32623.2495
$ws.Range("H139").Value2 = 7358.7856
$ws.Range("I139").Value2 = 3781.25
$ws.Range("K139").Value2 = 11343.75
$ws.Range("M139").Value2 = -6203.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value2 = 0
$ws.Range("J6").Value2 = 0
$ws.Range("L6").Value2 = 0
$ws.Range("N6").ClearContents() | Out-Null
$ws.Range("H16").Value2 = 0
$ws.Range("J16").Value2 = 0
$ws.Range("L16").Value2 = 0
$ws.Range("N16").ClearContents() | Out-Null
$ws.Range("H97").Value2 = 565.34375
$ws.Range("I97").Value2 = 642.3913
$ws.Range("K97").Value2 = 642.3913
$ws.Range("M97").Value2 = -146.3913
$ws.Range("H102").Value2 = 4056.6155
$ws.Range("I102").Value2 = 2031.9333
$ws.Range("K102").Value2 = 2031.9333
$ws.Range("M102").Value2 = -409.9332999999999
$ws.Range("H121").Value2 = 63998.5
$ws.Range("J121").Value2 = 63998.5
$ws.Range("L121").Value2 = 63998.5
$ws.Range("N121").Value2 = -67492.5
$ws.Range("H132").Value2 = 41902.61
$ws.Range("I132").Value2 = 5357.049
$ws.Range("K132").Value2 = 16071.147
$ws.Range("M132").Value2 = -13541.147
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H51").Value2 = 30000
$ws.Range("J51").Value2 = 30000
$ws.Range("L51").Value2 = 30000
$ws.Range("N51").Value2 = -30956
$ws.Range("H93").Value2 = 5194.952
$ws.Range("I93").Value2 = 4938.857
$ws.Range("J93").Value2 = 5707.143
$ws.Range("K93").Value2 = 4938.857
$ws.Range("L93").Value2 = 5707.143
$ws.Range("M93").Value2 = -3690.857
$ws.Range("N93").Value2 = -8203.143
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value2 = 1090.6666
$ws.Range("I113").Value2 = 974.7143
$ws.Range("K113").Value2 = 2924.1429
$ws.Range("M113").Value2 = -754.1428999999998
$ws.Range("H119").Value2 = 150000
$ws.Range("J119").Value2 = 150000
$ws.Range("L119").Value2 = 150000
$ws.Range("N119").Value2 = -159676
$ws.Range("H132").Value2 = 29963.285
$ws.Range("I132").Value2 = 1209.8438
$ws.Range("J132").Value2 = 336666.66
$ws.Range("K132").Value2 = 3629.5314
$ws.Range("L132").Value2 = 1009999.98
$ws.Range("M132").Value2 = -1099.5314
$ws.Range("N132").Value2 = -1015059.98
$ws.Range("H136").Value2 = 391056.3
$ws.Range("I136").Value2 = 388683.5
$ws.Range("K136").Value2 = 1166050.5
$ws.Range("M136").Value2 = -1163500.5
